# Update betting-odds values in row 4 to match the new FlashScore data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 3.05
$ws.Range("I4").Value = 2.12
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 2.62
$ws.Range("S4").Value = 1.32
$ws.Range("T4").Value = 3.2
$ws.Range("X4").Value = 17.5
$ws.Range("Y4").Value = 10.75
$ws.Range("Z4").Value = 40
$ws.Range("AA4").Value = 24
$ws.Range("AB4").Value = 28
$ws.Range("AE4").Value = 12.5
$ws.Range("AI4").Value = 11.75
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 15.5
$ws.Range("AM4").Value = 22
$ws.Range("AN4").Value = 5.1
$ws.Range("AO4").Value = 16
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 75
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 3
$ws.Range("AW4").Value = 4.15
$ws.Range("AX4").Value = 10.5
$ws.Range("AY4").Value = 16.5
$ws.Range("AZ4").Value = 37
$ws.Range("BA4").Value = 60
$ws.Range("BB4").Value = 175
